$d = $word.ActiveDocument

# Title (appears twice: main Heading1 and the bold line near the end) -
# wdReplaceAll (the final 2) takes care of both occurrences in one call.
$d.Content.Find.Execute("Play Blazing Mammoth Free | A Prehistoric Themed Slot Game", $true, $false, $false, $false, $false, $true, 1, $false, "Play Blazing Mammoth Free", 2)

# "What we like" bullet list
$d.Content.Find.Execute("Gameplay features", $true, $false, $false, $false, $false, $true, 1, $false, "High-quality graphics and sound effects", 2)

$d.Content.Find.Execute("Visual design", $true, $false, $false, $false, $false, $true, 1, $false, "Unique prehistoric theme", 2)

$d.Content.Find.Execute("Payouts and betting range", $true, $false, $false, $false, $false, $true, 1, $false, "Game multipliers and free spins", 2)

$d.Content.Find.Execute("Ease of use", $true, $false, $false, $false, $false, $true, 1, $false, "Easy-to-use interface for beginners", 2)

# "What we don't like" bullet list
$d.Content.Find.Execute("Lack of unique bonus features", $true, $false, $false, $false, $false, $true, 1, $false, "Limited betting range", 2)

$d.Content.Find.Execute("Limited betting range for high rollers", $true, $false, $false, $false, $false, $true, 1, $false, "Lack of innovative features", 2)

# Meta description (italic text near the end)
$d.Content.Find.Execute("Read our unbiased review of Blazing Mammoth, a 5-reel virtual slot machine with prehistoric theme. Play it for free and enjoy unique game features!", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Blazing Mammoth and play for free. Enjoy high-quality graphics and unique prehistoric theme.", 2)
